$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at the top of the date-ordered data block (rows 955-957),
# shifting the existing data (old rows 955-1024) down to rows 958-1027.
$ws.Range("A955:A957").EntireRow.Insert()

# Row 955: Fecha=2023-08-09 (45147), Variedad=Larga vida, Calidad=Primera
$ws.Cells.Item(955,1).Value = 2
$ws.Cells.Item(955,2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(955,3).Value = "Coquimbo"
$ws.Cells.Item(955,4).Value = "2023-08-09"
$ws.Cells.Item(955,5).Value = 4
$ws.Cells.Item(955,6).Value = 100112020
$ws.Cells.Item(955,7).Value = "Tomate"
$ws.Cells.Item(955,8).Value = "Larga vida"
$ws.Cells.Item(955,9).Value = "Primera"
$ws.Cells.Item(955,10).Value = 1000
$ws.Cells.Item(955,11).Value = 16000
$ws.Cells.Item(955,12).Value = 17000
$ws.Cells.Item(955,13).Value = 16500
$ws.Cells.Item(955,14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(955,15).Value = "Provincia de Limarí"
$ws.Cells.Item(955,16).Value = 917
$ws.Cells.Item(955,17).Value = 18
$ws.Cells.Item(955,18).Value = "Hortaliza"

# Row 956: Fecha=2023-08-09 (45147), Variedad=Larga vida, Calidad=Segunda
$ws.Cells.Item(956,1).Value = 2
$ws.Cells.Item(956,2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(956,3).Value = "Coquimbo"
$ws.Cells.Item(956,4).Value = "2023-08-09"
$ws.Cells.Item(956,5).Value = 4
$ws.Cells.Item(956,6).Value = 100112020
$ws.Cells.Item(956,7).Value = "Tomate"
$ws.Cells.Item(956,8).Value = "Larga vida"
$ws.Cells.Item(956,9).Value = "Segunda"
$ws.Cells.Item(956,10).Value = 1300
$ws.Cells.Item(956,11).Value = 13000
$ws.Cells.Item(956,12).Value = 14000
$ws.Cells.Item(956,13).Value = 13500
$ws.Cells.Item(956,14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(956,15).Value = "Provincia de Limarí"
$ws.Cells.Item(956,16).Value = 750
$ws.Cells.Item(956,17).Value = 18
$ws.Cells.Item(956,18).Value = "Hortaliza"

# Row 957: Fecha=2023-08-09 (45147), Variedad=Larga vida, Calidad=Tercera
$ws.Cells.Item(957,1).Value = 2
$ws.Cells.Item(957,2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(957,3).Value = "Coquimbo"
$ws.Cells.Item(957,4).Value = "2023-08-09"
$ws.Cells.Item(957,5).Value = 4
$ws.Cells.Item(957,6).Value = 100112020
$ws.Cells.Item(957,7).Value = "Tomate"
$ws.Cells.Item(957,8).Value = "Larga vida"
$ws.Cells.Item(957,9).Value = "Tercera"
$ws.Cells.Item(957,10).Value = 900
$ws.Cells.Item(957,11).Value = 10000
$ws.Cells.Item(957,12).Value = 11000
$ws.Cells.Item(957,13).Value = 10500
$ws.Cells.Item(957,14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(957,15).Value = "Provincia de Limarí"
$ws.Cells.Item(957,16).Value = 583
$ws.Cells.Item(957,17).Value = 18
$ws.Cells.Item(957,18).Value = "Hortaliza"
